$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace each customer's phone number with a masked placeholder value
$ws.Range("B2").Value = "55 00 0000-0000"
$ws.Range("B3").Value = "55 00 0000-0000"
$ws.Range("B4").Value = "55 00 0000-0000"
$ws.Range("B5").Value = "55 00 0000-0000"

# Make the phone number column formatting consistent (wrap text like B2 originally had)
$ws.Range("B2:B5").WrapText = $true

# Remove the now-unused trailing empty row (previously D6 held a style-only placeholder)
$ws.Rows("6:6").Delete()

# Update the active selection to reflect where the user left off editing
$ws.Range("F6").Select()
